$d = $word.ActiveDocument

# Step 1: Replace the original paragraph text with "Ответ 1" (literal Find/Replace keeps rPr, avoids xml:space="preserve").
$null = $d.Content.Find.Execute("На схему был добавлен конденсатор. Внутренний и внешний радиусы цилиндрического конденсатора изменены до 2 см и 2 см соответственно, также высота конденсатора изменена до 2 см, и значение относительной диэлектрической проницаемости изменено до 2. Внутренний и внешний радиусы цилиндрического конденсатора изменены до 0 см и 2 см соответственно, также высота конденсатора изменена до 2 см, и значение относительной диэлектрической проницаемости изменено до 2. Внутренний и внешний радиусы цилиндрического конденсатора изменены до 0 см и 0 см соответственно, также высота конденсатора изменена до 0 см, и значение относительной диэлектрической проницаемости изменено до 0. ", $false, $false, $false, $false, $false, $true, 1, $false, "Ответ 1", 2)

$orig = $d.Paragraphs(1)

# Step 2: Insert the "Вопрос №1" paragraph before it, and the "Контрольные вопросы:" heading before that.
$null = $orig.Range.InsertParagraphBefore()
$null = $d.Paragraphs(1).Range.InsertParagraphBefore()

# $orig is now paragraph 3 (index shifted after the two inserts above).
$p3 = $d.Paragraphs(3)

# Step 3: Insert "Вопрос №2" after it, then "Ответ 2" after that.
$null = $p3.Range.InsertParagraphAfter()
$null = $d.Paragraphs(4).Range.InsertParagraphAfter()

# --- Now configure each paragraph's formatting + text ---

# Paragraph 1: Heading "Контрольные вопросы:"
$p1 = $d.Paragraphs(1)
$p1.Style = "Heading 1"
$p1.Format.LineSpacing = 22.0
$p1.Format.FirstLineIndent = 42.75
$p1.Range.Text = "Контрольные вопросы:"
$p1TextRange = $d.Paragraphs(1).Range
$null = $p1TextRange.MoveEnd(1, -1)
$p1TextRange.Font.Bold = $true

# Paragraph 2: "Вопрос №1"
$p2 = $d.Paragraphs(2)
$p2.Format.Alignment = 0
$p2.Range.Text = "Вопрос №1"

# Paragraph 3: "Ответ 1" formatting (text already set via Find/Replace above)
$p3f = $d.Paragraphs(3)
$p3f.Format.SpaceAfter = 10.0

# Paragraph 4: "Вопрос №2"
$p4 = $d.Paragraphs(4)
$p4.Format.Alignment = 0
$p4.Range.Text = "Вопрос №2"

# Paragraph 5: "Ответ 2"
$p5 = $d.Paragraphs(5)
$p5.Format.SpaceAfter = 10.0
$p5.Range.Text = "Ответ 2"
